# Applies the "Updated cryptos list" price/volume refresh to sheet1.
# Rows 8/9 (Cardano/Dogecoin) swap position, row 51 (USDD -> Mantle) is
# replaced, and most other rows get refreshed Price/Volume(1h) text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.791.75'
$ws.Range("E2").Value = '  +0.63%  '
$ws.Range("D3").Value = '1.645.02'
$ws.Range("E3").Value = '  +0.11%  '
$ws.Range("E4").Value = '  +0.52%  '
$ws.Range("D5").Value = '''216.71'
$ws.Range("E5").Value = '  +0.36%  '
$ws.Range("E6").Value = '  -0.49%  '
$ws.Range("B8").Value = 'Dogecoin'
$ws.Range("C8").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D8").Value = '''0.0630'
$ws.Range("E8").Value = '  +0.57%  '
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").Value = '''0.251'
$ws.Range("E9").Value = '  -0.29%  '
$ws.Range("D10").Value = '''19.14'
$ws.Range("E10").Value = '  -0.36%  '
$ws.Range("E11").Value = '  +0.11%  '
$ws.Range("D12").Value = '1.868.93'
$ws.Range("E12").Value = '  -0.19%  '
$ws.Range("D13").Value = '1.635.85'
$ws.Range("E13").Value = '  -0.82%  '
$ws.Range("E14").Value = '  -0.97%  '
$ws.Range("E15").Value = '  -0.69%  '
$ws.Range("D16").Value = '''64.55'
$ws.Range("E16").Value = '  -2.08%  '
$ws.Range("D17").Value = '26.779.64'
$ws.Range("E17").Value = '  +0.37%  '
$ws.Range("E18").Value = '  -1.74%  '
$ws.Range("D19").Value = '''213.85'
$ws.Range("E19").Value = '  -2.05%  '
$ws.Range("E20").Value = '  +0.52%  '
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("D22").Value = '''2.47'
$ws.Range("E22").Value = '  +14.66%  '
$ws.Range("D23").Value = '''6.25'
$ws.Range("E23").Value = '  -0.92%  '
$ws.Range("D24").Value = '''9.35'
$ws.Range("E24").Value = '  -2.00%  '
$ws.Range("D25").Value = '''145.33'
$ws.Range("E25").Value = '  -0.62%  '
$ws.Range("E26").Value = '  +0.44%  '
$ws.Range("E27").Value = '  -1.44%  '
$ws.Range("D28").Value = '''7.11'
$ws.Range("E28").Value = '  -0.16%  '
$ws.Range("D29").Value = '''15.65'
$ws.Range("E29").Value = '  -1.27%  '
$ws.Range("D30").Value = '''0.0509'
$ws.Range("E30").Value = '  -1.87%  '
$ws.Range("E31").Value = '  +0.44%  '
$ws.Range("E32").Value = '  -1.84%  '
$ws.Range("E33").Value = '  -1.78%  '
$ws.Range("D34").Value = '1.296.74'
$ws.Range("E34").Value = '  +1.76%  '
$ws.Range("E35").Value = '  -0.31%  '
$ws.Range("E36").Value = '  +1.48%  '
$ws.Range("D37").Value = '''0.0174'
$ws.Range("E37").Value = '  -4.58%  '
$ws.Range("E38").Value = '  +0.90%  '
$ws.Range("D39").Value = '''0.825'
$ws.Range("E39").Value = '  -0.51%  '
$ws.Range("E40").Value = '  +0.52%  '
$ws.Range("D41").Value = '''0.810'
$ws.Range("E41").Value = '  +0.20%  '
$ws.Range("E42").Value = '  -0.20%  '
$ws.Range("E43").Value = '  -2.02%  '
$ws.Range("D44").Value = '1.795.51'
$ws.Range("D45").Value = '''61.77'
$ws.Range("E45").Value = '  +3.44%  '
$ws.Range("D46").Value = '''91.54'
$ws.Range("E46").Value = '  -1.74%  '
$ws.Range("E47").Value = '  +1.00%  '
$ws.Range("D48").Value = '''0.0524'
$ws.Range("E48").Value = '  +1.36%  '
$ws.Range("D49").Value = '''7.67'
$ws.Range("E49").Value = '  -1.95%  '
$ws.Range("D50").Value = '''0.0976'
$ws.Range("E50").Value = '  -0.17%  '
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").Value = '''0.408'
$ws.Range("E51").Value = '  +0.17%  '
